$d = $word.ActiveDocument

$replacements = @(
    @{old = "10÷4="; new = "11÷5="},
    @{old = "62÷4="; new = "28÷3="},
    @{old = "41÷3="; new = "68÷8="},
    @{old = "16÷6="; new = "21÷2="},
    @{old = "15÷2="; new = "98÷4="},
    @{old = "87÷8="; new = "64÷7="},
    @{old = "90÷4="; new = "88÷8="},
    @{old = "73÷4="; new = "11÷7="},
    @{old = "58÷3="; new = "82÷3="},
    @{old = "89÷8="; new = "28÷8="},
    @{old = "85÷4="; new = "93÷2="},
    @{old = "67÷3="; new = "22÷2="},
    @{old = "25÷2="; new = "20÷6="},
    @{old = "68÷3="; new = "17÷9="},
    @{old = "99÷8="; new = "49÷4="},
    @{old = "53÷2="; new = "48÷3="},
    @{old = "29÷8="; new = "67÷6="},
    @{old = "75÷7="; new = "61÷2="},
    @{old = "66÷3="; new = "19÷5="},
    @{old = "82÷4="; new = "11÷3="},
    @{old = "36÷8="; new = "98÷9="},
    @{old = "18÷4="; new = "60÷6="},
    @{old = "73÷5="; new = "47÷6="},
    @{old = "86÷7="; new = "99÷5="},
    @{old = "44÷4="; new = "23÷7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
